$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format numeric-looking text cells as Text so Excel keeps them as literal strings
$textCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D32", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D46", "D47", "D48", "D49", "D50")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '63.165.82'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '3.183.62'
$ws.Range("E3").Value = '  -3.92%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '592.75'
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("D6").Value = '135.42'
$ws.Range("E6").Value = '  -4.33%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.180.97'
$ws.Range("E8").Value = '  -3.96%  '
$ws.Range("D9").Value = '0.515'
$ws.Range("E9").Value = '  -0.75%  '
$ws.Range("D10").Value = '0.141'
$ws.Range("E10").Value = '  -5.89%  '
$ws.Range("D11").Value = '5.23'
$ws.Range("E11").Value = '  -5.36%  '
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  -3.08%  '
$ws.Range("D13").Value = '0.0000237'
$ws.Range("E13").Value = '  -4.36%  '
$ws.Range("D14").Value = '34.74'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '3.710.06'
$ws.Range("E15").Value = '  -3.86%  '
$ws.Range("E16").Value = '  -1.01%  '
$ws.Range("D17").Value = '3.186.24'
$ws.Range("E17").Value = '  -3.84%  '
$ws.Range("D18").Value = '63.097.38'
$ws.Range("E18").Value = '  -1.04%  '
$ws.Range("D19").Value = '6.58'
$ws.Range("E19").Value = '  -4.05%  '
$ws.Range("D20").Value = '462.57'
$ws.Range("E20").Value = '  -3.77%  '
$ws.Range("D21").Value = '14.03'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").Value = '0.698'
$ws.Range("E22").Value = '  -5.61%  '
$ws.Range("D23").Value = '7.65'
$ws.Range("E23").Value = '  -4.32%  '
$ws.Range("D24").Value = '13.39'
$ws.Range("E24").Value = '  -4.38%  '
$ws.Range("D25").Value = '82.83'
$ws.Range("E25").Value = '  -2.74%  '
$ws.Range("B27").Value = 'FirstDigitalUSD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.68'
$ws.Range("E28").Value = '  -3.61%  '
$ws.Range("D29").Value = '7.71'
$ws.Range("E29").Value = '  -6.14%  '
$ws.Range("D30").Value = '6.77'
$ws.Range("E30").Value = '  -5.69%  '
$ws.Range("E31").Value = '  -5.23%  '
$ws.Range("D32").Value = '27.25'
$ws.Range("E32").Value = '  -5.60%  '
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("E34").Value = '  -5.51%  '
$ws.Range("E35").Value = '  -6.33%  '
$ws.Range("D36").Value = '5.82'
$ws.Range("E36").Value = '  -3.93%  '
$ws.Range("D37").Value = '51.32'
$ws.Range("E37").Value = '  -2.25%  '
$ws.Range("D38").Value = '0.0₃0709'
$ws.Range("E38").Value = '  -5.18%  '
$ws.Range("D39").Value = '0.0389'
$ws.Range("E39").Value = '  -2.77%  '
$ws.Range("D40").Value = '406.81'
$ws.Range("E40").Value = '  -6.04%  '
$ws.Range("D41").Value = '8.09'
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("D42").Value = '2.65'
$ws.Range("E42").Value = '  -3.14%  '
$ws.Range("D43").Value = '0.112'
$ws.Range("E43").Value = '  -5.80%  '
$ws.Range("D44").Value = '2.809.22'
$ws.Range("E44").Value = '  -9.93%  '
$ws.Range("E45").Value = '  -5.38%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '0.999'
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '2.12'
$ws.Range("E47").Value = '  -4.92%  '
$ws.Range("D48").Value = '35.15'
$ws.Range("E48").Value = '  -4.52%  '
$ws.Range("D49").Value = '124.62'
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("D50").Value = '25.27'
$ws.Range("E50").Value = '  -3.89%  '
$ws.Range("E51").Value = '  -1.80%  '

# Restore default style for cells we temporarily formatted as text
foreach ($c in $textCells) {
    $ws.Range($c).Style = "Normal"
}